# Add the new "Enero" rows (17, 19, 20, 22 -- note rows 18 and 21 are left
# blank, matching the source diff) to the "Horas2018" sheet, then update the
# selection to I8 to match the post-edit workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Horas2018")

# New rows of data: Mes (B), Dia (C), Horas (D)
$ws.Cells.Item(17, 2).Value = "Enero"
$ws.Cells.Item(17, 3).Value = 4
$ws.Cells.Item(17, 4).Value = 3

$ws.Cells.Item(19, 2).Value = "Enero"
$ws.Cells.Item(19, 3).Value = 7
$ws.Cells.Item(19, 4).Value = 4

$ws.Cells.Item(20, 2).Value = "Enero"
$ws.Cells.Item(20, 3).Value = 8
$ws.Cells.Item(20, 4).Value = 3

$ws.Cells.Item(22, 2).Value = "Enero"
$ws.Cells.Item(22, 3).Value = 11
$ws.Cells.Item(22, 4).Value = 5

# Copy the formatting from the last existing data row (16) onto each of the
# newly populated rows only (rows 18 and 21 stay completely untouched/empty,
# matching the source data).
$ws.Range("B16:D16").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)
$ws.Range("B19:D19").PasteSpecial(-4122)
$ws.Range("B20:D20").PasteSpecial(-4122)
$ws.Range("B22:D22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Extend the summary formulas to cover the new rows.
$ws.Range("H3").Formula = "=COUNT(C4:C22)"
$ws.Range("H4").Formula = "=SUM(D4:D22)"

$ws.Calculate()

# Match the saved selection from the edited workbook.
$ws.Range("I8").Select()
